{"js": "// R\u00e9sum\u00e9 update:\n//  1. The career-history row for \"\uc8fc\uc2dd\ud68c\uc0ac \uc560\ub4dc\uc5c5\" had an end-date of\n//     \"\ud604  \uc7ac\" (Korean for \"present\" / still employed). The employee has\n//     since left that position, so the end date is changed to a fixed\n//     date: \"2024-12\". The preceding \"2018.12 ~ \" text is left untouched.\n//  2. Drop the stray \"_GoBack\" bookmark that Word leaves behind after an\n//     editing session (harmless leftover cleaned up on save).\n\nconst body = context.document.body;\n\n// Locate the exact \"\ud604  \uc7ac\" (\ud604, space, space, \uc7ac) text \u2014 this string is\n// unique in the document (the other \"\ud604 \uc7ac\" occurrence, in the project\n// table's \"2021.06 ~ \ud604   \uc7ac\" row, uses different spacing and is left\n// alone).\nconst results = body.search(\"\ud604  \uc7ac\", { matchCase: true, matchWildcards: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(`Expected exactly one match for \"\ud604  \uc7ac\", found ${results.items.length}`);\n}\n\nresults.items[0].insertText(\"2024-12\", Word.InsertLocation.replace);\n\n// Remove the leftover \"_GoBack\" bookmark (no-op if it is not present).\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# R\u00e9sum\u00e9 update:\n#  1. The career-history row for \"\uc8fc\uc2dd\ud68c\uc0ac \uc560\ub4dc\uc5c5\" had an end-date of\n#     \"\ud604  \uc7ac\" (Korean for \"present\" / still employed). The employee has\n#     since left that position, so the end date is changed to a fixed\n#     date: \"2024-12\". The preceding \"2018.12 ~ \" text is left untouched.\n#  2. Drop the stray \"_GoBack\" bookmark that Word leaves behind after an\n#     editing session (harmless leftover cleaned up on save).\n\n$d = $word.ActiveDocument\n\n# Locate the exact \"\ud604  \uc7ac\" (\ud604, space, space, \uc7ac) text \u2014 this string is\n# unique in the document (the other \"\ud604 \uc7ac\" occurrence, in the project\n# table's \"2021.06 ~ \ud604   \uc7ac\" row, uses different spacing and is left\n# alone).\n$rng = $d.Content\n$found = $rng.Find.Execute(\"\ud604  \uc7ac\")\nif ($found) {\n    $rng.Text = \"2024-12\"\n} else {\n    throw \"Could not find the '\ud604  \uc7ac' text to update\"\n}\n\n# Remove the leftover \"_GoBack\" bookmark, if present.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
